$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1650
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H70").Value = 1066.3334
$ws.Range("I70").Value = 1099.5
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 3298.5
$ws.Range("L70").Value = 3000
$ws.Range("M70").Value = -3028.5
$ws.Range("N70").Value = -3540
$ws.Range("H73").Value = 1066.3334
$ws.Range("I73").Value = 1099.5
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 3298.5
$ws.Range("L73").Value = 3000
$ws.Range("M73").Value = -2362.5
$ws.Range("N73").Value = -4872
$ws.Range("H132").Value = 1201.7931
$ws.Range("I132").Value = 1201.7931
$ws.Range("K132").Value = 3605.379300000001
$ws.Range("M132").Value = -1075.379300000001
$ws.Range("H134").Value = 124800
$ws.Range("J134").Value = 124800
$ws.Range("L134").Value = 124800
$ws.Range("N134").Value = -134940
$ws.Range("H135").Value = 1997.375
$ws.Range("I135").Value = 1663.1666
$ws.Range("K135").Value = 14968.4994
$ws.Range("M135").Value = -12433.4994
$ws.Range("H138").Value = 4346.8335
$ws.Range("I138").Value = 3904.9412
$ws.Range("J138").Value = 4742.2104
$ws.Range("K138").Value = 11714.8236
$ws.Range("L138").Value = 14226.6312
$ws.Range("M138").Value = -6574.8236
$ws.Range("N138").Value = -24506.6312
$ws.Range("H141").Value = 8332.833000000001
$ws.Range("I141").Value = 7999.3335
$ws.Range("K141").Value = 23998.0005
$ws.Range("M141").Value = -18818.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10762.4
$ws.Range("I32").Value = 10459.982
$ws.Range("K32").Value = 10459.982
$ws.Range("M32").Value = -10172.982
$ws.Range("H45").Value = 3347.8333
$ws.Range("I45").Value = 2724.6667
$ws.Range("J45").Value = 3971
$ws.Range("K45").Value = 2724.6667
$ws.Range("L45").Value = 3971
$ws.Range("M45").Value = -2347.6667
$ws.Range("N45").Value = -4725
$ws.Range("H135").Value = 52084.5
$ws.Range("J135").Value = 52084.5
$ws.Range("L135").Value = 52084.5
$ws.Range("N135").Value = -62224.5
$ws.Range("H138").Value = 99893.5
$ws.Range("J138").Value = 99893.5
$ws.Range("L138").Value = 99893.5
$ws.Range("N138").Value = -110173.5
$ws.Range("H139").Value = 121715
$ws.Range("J139").Value = 121715
$ws.Range("L139").Value = 121715
$ws.Range("N139").Value = -131995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H81").Value = 74900
$ws.Range("J81").Value = 74900
$ws.Range("L81").Value = 74900
$ws.Range("N81").Value = -77022
$ws.Range("H84").Value = 74900
$ws.Range("J84").Value = 74900
$ws.Range("L84").Value = 224700
$ws.Range("N84").Value = -235308
$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492
$ws.Range("H105").Value = 7572
$ws.Range("I105").Value = 7572
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 7572
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -5825
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1081.7142
$ws.Range("I107").Value = 928.8333
$ws.Range("K107").Value = 928.8333
$ws.Range("M107").Value = 991.1667
$ws.Range("H134").Value = 3342.0908
$ws.Range("I134").Value = 3276.4
$ws.Range("K134").Value = 9829.200000000001
$ws.Range("M134").Value = -7294.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 719.3
$ws.Range("I22").Value = 689.25
$ws.Range("J22").Value = 839.5
$ws.Range("K22").Value = 689.25
$ws.Range("L22").Value = 839.5
$ws.Range("M22").Value = -339.25
$ws.Range("N22").Value = -1539.5
$ws.Range("H86").Value = 22811.684
$ws.Range("I86").Value = 10284.182
$ws.Range("J86").Value = 40037
$ws.Range("K86").Value = 10284.182
$ws.Range("L86").Value = 40037
$ws.Range("M86").Value = -9161.182000000001
$ws.Range("N86").Value = -42283
$ws.Range("H89").Value = 22811.684
$ws.Range("I89").Value = 10284.182
$ws.Range("J89").Value = 40037
$ws.Range("K89").Value = 51420.91
$ws.Range("L89").Value = 200185
$ws.Range("M89").Value = -45804.91
$ws.Range("N89").Value = -211417
$ws.Range("H99").Value = 9111.25
$ws.Range("I99").Value = 8778.200000000001
$ws.Range("K99").Value = 8778.200000000001
$ws.Range("M99").Value = -7280.200000000001
$ws.Range("H126").Value = 9111.25
$ws.Range("I126").Value = 8778.200000000001
$ws.Range("K126").Value = 26334.6
$ws.Range("M126").Value = -23864.6
$ws.Range("H134").Value = 3573.1428
$ws.Range("I134").Value = 3252
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 9756
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -7221
$ws.Range("N134").Value = -21570
$ws.Range("H141").Value = 58686.8
$ws.Range("J141").Value = 58686.8
$ws.Range("L141").Value = 58686.8
$ws.Range("N141").Value = -69046.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 220000
$ws.Range("I128").Value = 220000
$ws.Range("K128").Value = 660000
$ws.Range("M128").Value = -655020
$ws.Range("H129").Value = 1438.2
$ws.Range("I129").Value = 1113.8334
$ws.Range("J129").Value = 1924.75
$ws.Range("K129").Value = 3341.5002
$ws.Range("L129").Value = 5774.25
$ws.Range("M129").Value = 1658.4998
$ws.Range("N129").Value = -15774.25
$ws.Range("H131").Value = 1871.7858
$ws.Range("I131").Value = 1521.75
$ws.Range("J131").Value = 2011.8
$ws.Range("K131").Value = 4565.25
$ws.Range("L131").Value = 6035.4
$ws.Range("M131").Value = 474.75
$ws.Range("N131").Value = -16115.4
$ws.Range("H137").Value = 4480.6665
$ws.Range("I137").Value = 999.5
$ws.Range("J137").Value = 5176.9
$ws.Range("K137").Value = 2998.5
$ws.Range("L137").Value = 15530.7
$ws.Range("M137").Value = 2101.5
$ws.Range("N137").Value = -25730.7
$ws.Range("H140").Value = 1233.7693
$ws.Range("I140").Value = 1233.7693
$ws.Range("K140").Value = 3701.3079
$ws.Range("M140").Value = 1478.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3571.7273
$ws.Range("I122").Value = 2898.5715
$ws.Range("J122").Value = 4749.75
$ws.Range("K122").Value = 8695.7145
$ws.Range("L122").Value = 14249.25
$ws.Range("M122").Value = -6245.7145
$ws.Range("N122").Value = -19149.25
$ws.Range("H132").Value = 4373.8
$ws.Range("I132").Value = 4445
$ws.Range("J132").Value = 4241.5713
$ws.Range("K132").Value = 13335
$ws.Range("L132").Value = 12724.7139
$ws.Range("M132").Value = -10805
$ws.Range("N132").Value = -17784.7139
$ws.Range("H139").Value = 75673.125
$ws.Range("J139").Value = 75673.125
$ws.Range("L139").Value = 75673.125
$ws.Range("N139").Value = -85953.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3004
$ws.Range("I40").Value = 3004
$ws.Range("K40").Value = 3004
$ws.Range("M40").Value = -2868
$ws.Range("H55").Value = 385.875
$ws.Range("J55").Value = 386
$ws.Range("L55").Value = 386
$ws.Range("N55").Value = -732
$ws.Range("H132").Value = 5784.375
$ws.Range("I132").Value = 5712.8335
$ws.Range("K132").Value = 17138.5005
$ws.Range("M132").Value = -14608.5005
$ws.Range("H136").Value = 5203.1055
$ws.Range("I136").Value = 5561.857
$ws.Range("J136").Value = 4198.6
$ws.Range("K136").Value = 16685.571
$ws.Range("L136").Value = 12595.8
$ws.Range("M136").Value = -14135.571
$ws.Range("N136").Value = -17695.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H132").Value = 2636.2
$ws.Range("I132").Value = 1939
$ws.Range("K132").Value = 5817
$ws.Range("M132").Value = -3287
